$wb = $excel.ActiveWorkbook

# Worksheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 342
$ws.Range("F7").Value = 5958
$ws.Range("F9").Value = 4012
$ws.Range("F10").Value = 219
$ws.Range("F13").Value = 132
$ws.Range("F14").Value = 139
$ws.Range("F15").Value = 3999
$ws.Range("F16").Value = 14
$ws.Range("F19").Value = 5601
$ws.Range("F20").Value = 449
$ws.Range("F21").Value = 2206
$ws.Range("G21").Value = 70
$ws.Range("F23").Value = 397
$ws.Range("F24").Value = 8396
$ws.Range("F26").Value = 50
$ws.Range("F27").Value = 2242
$ws.Range("F28").Value = 2279
$ws.Range("F31").Value = 1865
$ws.Range("F32").Value = 37
$ws.Range("F33").Value = 296
$ws.Range("F35").Value = 13
$ws.Range("F37").Value = 27
$ws.Range("F38").Value = 19
$ws.Range("F39").Value = 1202
$ws.Range("F40").Value = 1195
$ws.Range("F41").Value = 50
$ws.Range("F42").Value = 73
$ws.Range("F43").Value = 196
$ws.Range("F44").Value = 1394
$ws.Range("F45").Value = 2245
$ws.Range("F48").Value = 1230

# Worksheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1
$ws.Range("F13").Value = 135
$ws.Range("F22").Value = 24

# Worksheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 828
$ws.Range("F4").Value = 80

# Worksheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 828
$ws.Range("F6").Value = 342
$ws.Range("F7").Value = 5958
$ws.Range("F9").Value = 4012
$ws.Range("F10").Value = 219
$ws.Range("F13").Value = 132
$ws.Range("F16").Value = 139
$ws.Range("F17").Value = 3999
$ws.Range("F18").Value = 14
$ws.Range("F21").Value = 5601
$ws.Range("F22").Value = 449
$ws.Range("F23").Value = 2206
$ws.Range("G23").Value = 70
$ws.Range("F25").Value = 397
$ws.Range("F26").Value = 8396
$ws.Range("F27").Value = 135
$ws.Range("F29").Value = 2242
$ws.Range("F30").Value = 2279
$ws.Range("F33").Value = 1865
$ws.Range("F34").Value = 37
$ws.Range("F35").Value = 296
$ws.Range("F36").Value = 13
$ws.Range("F38").Value = 27
$ws.Range("F39").Value = 19
$ws.Range("F40").Value = 1202
$ws.Range("F41").Value = 1195
$ws.Range("F42").Value = 50
$ws.Range("F43").Value = 73
$ws.Range("F44").Value = 196
$ws.Range("F45").Value = 1394
$ws.Range("F46").Value = 2245
$ws.Range("F49").Value = 1230
